$d = $word.ActiveDocument

# Step 1: duplicate the "Title 5" paragraph and insert the copy right after
# the empty paragraph that follows it (this produces a brand-new paragraph
# with identical run formatting, rather than a bare text insertion).
$pTitle5 = $d.Paragraphs(13)
$srcTitle5 = $d.Range($pTitle5.Range.Start, $pTitle5.Range.End)
$pAfterTitle5 = $d.Paragraphs(14)
$insertPoint1 = $d.Range($pAfterTitle5.Range.End, $pAfterTitle5.Range.End)
$insertPoint1.FormattedText = $srcTitle5.FormattedText

# Step 2: change the text of the newly inserted paragraph from "Title 5" to
# "Title 6".
$pNewTitle = $d.Paragraphs(15)
$pNewTitle.Range.Find.Execute("Title 5", $false, $false, $false, $false, $false, $true, 1, $false, "Title 6", 2)

# Step 3: duplicate the new "Title 6" paragraph to create a following
# paragraph with the same formatting, then strip its text so only an empty
# paragraph remains (matching the style of the other blank-line paragraphs).
$pTitle6 = $d.Paragraphs(15)
$srcTitle6 = $d.Range($pTitle6.Range.Start, $pTitle6.Range.End)
$insertPoint2 = $d.Range($pTitle6.Range.End, $pTitle6.Range.End)
$insertPoint2.FormattedText = $srcTitle6.FormattedText

$pNewEmpty = $d.Paragraphs(16)
$clearRange = $d.Range($pNewEmpty.Range.Start, $pNewEmpty.Range.End - 1)
$clearRange.Text = ""

# Step 4: move the "_GoBack" bookmark from its old location (the empty
# paragraph right after "Title 4") to the very last paragraph of the
# document.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
$bookmarkRange = $lastParagraph.Range.Duplicate
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
